$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 350, pushing the existing rows 350-378 down to 352-380.
$ws.Rows.Item(350).Resize(2).Insert()

# Row 350: new weekly entry (Primera)
$ws.Cells.Item(350, 1).Value = 3
$ws.Cells.Item(350, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(350, 3).Value = "Coquimbo"
$ws.Cells.Item(350, 4).Value = 44461
$ws.Cells.Item(350, 5).Value = 5
$ws.Cells.Item(350, 6).Value = 100114014
$ws.Cells.Item(350, 7).Value = "Betarraga"
$ws.Cells.Item(350, 8).Value = "Sin especificar"
$ws.Cells.Item(350, 9).Value = "Primera"
$ws.Cells.Item(350, 10).Value = 3000
$ws.Cells.Item(350, 11).Value = 500
$ws.Cells.Item(350, 12).Value = 550
$ws.Cells.Item(350, 13).Value = 527
$ws.Cells.Item(350, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(350, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(350, 16).Value = 132
$ws.Cells.Item(350, 17).Value = 4
$ws.Cells.Item(350, 18).Value = "Hortaliza"

# Row 351: new weekly entry (Segunda)
$ws.Cells.Item(351, 1).Value = 3
$ws.Cells.Item(351, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(351, 3).Value = "Coquimbo"
$ws.Cells.Item(351, 4).Value = 44461
$ws.Cells.Item(351, 5).Value = 5
$ws.Cells.Item(351, 6).Value = 100114014
$ws.Cells.Item(351, 7).Value = "Betarraga"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "Segunda"
$ws.Cells.Item(351, 10).Value = 1500
$ws.Cells.Item(351, 11).Value = 400
$ws.Cells.Item(351, 12).Value = 400
$ws.Cells.Item(351, 13).Value = 400
$ws.Cells.Item(351, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(351, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(351, 16).Value = 100
$ws.Cells.Item(351, 17).Value = 4
$ws.Cells.Item(351, 18).Value = "Hortaliza"
